# Plots.xlsx: add a new DataCombined entry "AciclovirPVB" (row 2, column A)
# and make the DataCombined sheet the active/selected sheet (it previously
# wasn't; "plotGrids" was active/tabSelected before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataCombined")

# New data row: name = AciclovirPVB (other columns on row 2 stay blank,
# same as the existing blank template row).
$ws.Range("A2").Value = "AciclovirPVB"

# Make DataCombined the active sheet with A2 selected.
$ws.Activate()
$ws.Range("A2").Select()
